$wb = $excel.ActiveWorkbook

# Fix trailing space in the first sheet's name:
# "mst_patient_registration " -> "mst_patient_registration"
$wsRegistration = $wb.Worksheets.Item("mst_patient_registration ")
$wsRegistration.Name = "mst_patient_registration"

# Switch the active/selected tab from "mst_features" (last sheet) to the
# first sheet ("mst_patient_registration"), which also clears the old
# tabSelected flag on "mst_features" and the firstSheet attribute.
$wsRegistration.Activate()
